$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$dx = 294928 / 12700.0
$dy = -1567296 / 12700.0

# Move every shape except the empty Title placeholder by (dx, dy) EMU.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -ne "Title 26") {
        $sh.Left = $sh.Left + $dx
        $sh.Top = $sh.Top + $dy
    }
}

# Remove the empty title placeholder shape.
$s.Shapes.Item("Title 26").Delete()

Write-Host "done"
